# Rename the original sheet and add a second "Info" sheet, then populate it,
# matching the target diff (renames Sheet1 -> Product, adds Info sheet).

$wb = $excel.ActiveWorkbook

# --- Sheet1 -> Product ---------------------------------------------------
$wsProduct = $wb.Worksheets.Item(1)
$wsProduct.Name = "Product"

# Update Product sheet selection: select A1:F8, no tab-select override needed
# (tabSelected moves automatically once Info becomes active).
[void]$wsProduct.Range("A1:F8").Select()

# --- Add the new Info sheet, placed right after Product ------------------
$wsInfo = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsProduct)
$wsInfo.Name = "Info"

# Row 1 - headers
$wsInfo.Range("A1").Value = "em"
$wsInfo.Range("D1").Value = "Width"
$wsInfo.Range("D1:F1").HorizontalAlignment = -4108
[void]$wsInfo.Range("D1:F1").Merge()

# Row 2
$wsInfo.Range("A2").Value = "rem"
$wsInfo.Range("B2").Value = 18
$wsInfo.Range("D2").Value = 960
$wsInfo.Range("E2").Formula = '=$D$2*0.35'
$wsInfo.Range("G2").Value = 2
$wsInfo.Range("H2").Value = "rem"
$wsInfo.Range("I2").Formula = '=G2*$B$2'

# Row 3
$wsInfo.Range("E3").Formula = '=$D$2*0.35'
$wsInfo.Range("G3").Formula = '=ROUND(1-SUM(I2:I4)/E2,2)*100'
$wsInfo.Range("H3").Value = "%"

# Row 4
$wsInfo.Range("G4").Value = 2
$wsInfo.Range("H4").Value = "rem"
$wsInfo.Range("I4").Formula = '=G4*$B$2'

# Row 5
$wsInfo.Range("D5").Formula = '=SUM(E2:E4)'

# Row 6
$wsInfo.Range("E6").Formula = '=(D2-D5)'

# Row 8
$wsInfo.Range("E8").Formula = '=SUM(E2:E7)'

# Info becomes the active sheet / tab, with D3 selected.
[void]$wsInfo.Activate()
[void]$wsInfo.Range("D3").Select()
